# RCL-B-CDH1 Data Budget.xlsx - "What did i do?"
#
# 1. Bumped the FPGA bit-depth assumption on the "Data Budget" sheet from
#    24*8 (192 bits) to 32*8 (256 bits) for the first mission segment (C5),
#    and re-pointed the other two segments (C14, C23) at that same cell
#    instead of re-typing the literal 24*8 formula, so they track C5.
# 2. Left the Data Budget sheet active/selected on C4 (where the work was
#    happening) and tidied the Assumptions sheet selection down to a single
#    cell (C2) instead of the C2:C6 block.
#
# All the downstream totals (E5/E6/E8, E14/E15/E17, E23/E24/E26, E28, E30)
# and the mirrored figures on "For Pie" (C2,C3,C4,C6) are plain formulas
# that recalc from the above - no direct edits needed there.

$wb = $excel.ActiveWorkbook

$assumptions = $wb.Worksheets.Item("Assumptions")
$dataBudget  = $wb.Worksheets.Item("Data Budget")

# Assumptions sheet: shrink the lingering selection block down to one cell.
[void]$assumptions.Range("C2").Select()

# Data Budget sheet: the actual edit - bump the FPGA byte count and make
# the other two segments reference the first instead of duplicating it.
$dataBudget.Range("C5").Formula = "=32*8"
$dataBudget.Range("C14").Formula = "=C5"
$dataBudget.Range("C23").Formula = "=C5"

# Leave the workbook focused on the Data Budget sheet, selection on C4.
[void]$dataBudget.Activate()
[void]$dataBudget.Range("C4").Select()
